# Daily attendance processing - 2025-12-23 11:00:14
# Rotate the comma-separated "Recorded By" values in column G so that the
# last entry in the list is moved to the front (e.g. "System, dnasr281@gmail.com"
# becomes "dnasr281@gmail.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Count -gt 1) {
        $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
        $cell.Value2 = [string]::Join(", ", $rotated)
    }
}
